$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Identificadores de CU (B) primero
$ws.Range("B27").Value = "CU - 23"
$ws.Range("B28").Value = "CU - 24"

# Alias (D)
$ws.Range("D27").Value = "Consultar profesores"
$ws.Range("D28").Value = "Consultar clientes"

# Descripcion de Caso de Uso (C)
$ws.Range("C27").Value = "El director puede  consultar todos los profesores en la institución."
$ws.Range("C28").Value = "El director puede consultar todos los clientes registrados."

# Estado (E)
$ws.Range("E27").Value = "vacio"
$ws.Range("E28").Value = "vacio"

# Esfuerzo (F)
$ws.Range("F27").Value = 0
$ws.Range("F28").Value = 0

# Incremento (G)
$ws.Range("G27").Value = 0
$ws.Range("G28").Value = 0

# Prioridad (H)
$ws.Range("H27").Value = 1
$ws.Range("H28").Value = 1

$ws.Range("C28").Select()
